$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 331; this pushes the existing rows 331-338 down to
# 332-339 (dimension grows from A1:R338 to A1:R339).
$ws.Rows("331:331").Insert()

# Populate the newly inserted row 331 with the new weekly price record.
$ws.Cells.Item(331, 1).Value = 6
$ws.Cells.Item(331, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(331, 3).Value = "Metropolitana"
$ws.Cells.Item(331, 4).Value = 44448
$ws.Cells.Item(331, 5).Value = 13
$ws.Cells.Item(331, 6).Value = 100112044
$ws.Cells.Item(331, 7).Value = "Perejil"
$ws.Cells.Item(331, 8).Value = "Sin especificar"
$ws.Cells.Item(331, 9).Value = "Primera"
$ws.Cells.Item(331, 10).Value = 260
$ws.Cells.Item(331, 11).Value = 8000
$ws.Cells.Item(331, 12).Value = 9000
$ws.Cells.Item(331, 13).Value = 8423
$ws.Cells.Item(331, 14).Value = "`$/docena de atados"
$ws.Cells.Item(331, 15).Value = "Región Metropolitana"
$ws.Cells.Item(331, 16).Value = 2808
$ws.Cells.Item(331, 17).Value = 3
$ws.Cells.Item(331, 18).Value = "Hortaliza"
